# Add a new worksheet "tblStudyManagementTools1" that carries the same
# tabular data as the original sheet, then register a workbook-level
# defined name "tblStudyManagementTools" pointing at its A1:D77 range.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New worksheet placed right after the original sheet.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "tblStudyManagementTools1"

$ws2.Cells.Item(1, 1).Value = 'ID'
$ws2.Cells.Item(1, 2).Value = 'fkStudyID'
$ws2.Cells.Item(1, 3).Value = 'MPManagementTool'
$ws2.Cells.Item(1, 4).Value = 'MPAlternativesEvaluated'
$ws2.Cells.Item(2, 1).Value = 1
$ws2.Cells.Item(2, 2).Value = 1
$ws2.Cells.Item(2, 3).Value = 'Catch Limit'
$ws2.Cells.Item(3, 1).Value = 2
$ws2.Cells.Item(3, 2).Value = 2
$ws2.Cells.Item(3, 3).Value = 'Catch Limit'
$ws2.Cells.Item(4, 1).Value = 3
$ws2.Cells.Item(4, 2).Value = 3
$ws2.Cells.Item(4, 3).Value = 'Catch Limit'
$ws2.Cells.Item(5, 1).Value = 4
$ws2.Cells.Item(5, 2).Value = 4
$ws2.Cells.Item(5, 3).Value = 'Catch Limit'
$ws2.Cells.Item(6, 1).Value = 5
$ws2.Cells.Item(6, 2).Value = 5
$ws2.Cells.Item(6, 3).Value = 'Catch Limit'
$ws2.Cells.Item(7, 1).Value = 6
$ws2.Cells.Item(7, 2).Value = 6
$ws2.Cells.Item(7, 3).Value = 'Catch Limit'
$ws2.Cells.Item(8, 1).Value = 7
$ws2.Cells.Item(8, 2).Value = 6
$ws2.Cells.Item(8, 3).Value = 'Share Allocation'
$ws2.Cells.Item(9, 1).Value = 8
$ws2.Cells.Item(9, 2).Value = 7
$ws2.Cells.Item(9, 3).Value = 'Effort Limit'
$ws2.Cells.Item(10, 1).Value = 13
$ws2.Cells.Item(10, 2).Value = 11
$ws2.Cells.Item(10, 3).Value = 'Effort Limit'
$ws2.Cells.Item(10, 4).Value = 'Four effort levels. Low, Moderate, High, and Very High effort equivalent to 80t, 110t, 140t, and 170t harvests.'
$ws2.Cells.Item(11, 1).Value = 14
$ws2.Cells.Item(11, 2).Value = 11
$ws2.Cells.Item(11, 3).Value = 'Closure'
$ws2.Cells.Item(11, 4).Value = 'Two seasonal closures. Seasonal closure for 2 months or no closure.'
$ws2.Cells.Item(12, 1).Value = 15
$ws2.Cells.Item(12, 2).Value = 11
$ws2.Cells.Item(12, 3).Value = 'Closure'
$ws2.Cells.Item(12, 4).Value = 'Two spatial closures. Spatial closure of 10 nautical miles around island or no closure.'
$ws2.Cells.Item(13, 1).Value = 16
$ws2.Cells.Item(13, 2).Value = 11
$ws2.Cells.Item(13, 3).Value = 'Size Limit'
$ws2.Cells.Item(13, 4).Value = 'Three minimum size limits. 35, 38, and 40 cm TL.'
$ws2.Cells.Item(14, 1).Value = 17
$ws2.Cells.Item(14, 2).Value = 12
$ws2.Cells.Item(14, 3).Value = 'Catch Limit'
$ws2.Cells.Item(14, 4).Value = 'Levels of adjustement to the catch limit to account for uncertainty'
$ws2.Cells.Item(15, 1).Value = 18
$ws2.Cells.Item(15, 2).Value = 13
$ws2.Cells.Item(15, 3).Value = 'Closure'
$ws2.Cells.Item(15, 4).Value = 'Three levels of area closure. Reef perimeter closed to fishing. 16%, 32%, 50% closed.'
$ws2.Cells.Item(16, 1).Value = 19
$ws2.Cells.Item(16, 2).Value = 13
$ws2.Cells.Item(16, 3).Value = 'Effort Limit'
$ws2.Cells.Item(16, 4).Value = 'Three levels of annual fishing effort. 0.5, 1, and 1.5 times 1996 effort.'
$ws2.Cells.Item(17, 1).Value = 20
$ws2.Cells.Item(17, 2).Value = 14
$ws2.Cells.Item(17, 3).Value = 'Catch Limit'
$ws2.Cells.Item(17, 4).Value = 'Harvest Control Rule sets target F, fishery mortality level.  HCR implementation method not specified.'
$ws2.Cells.Item(18, 1).Value = 21
$ws2.Cells.Item(18, 2).Value = 15
$ws2.Cells.Item(18, 3).Value = 'Catch Limit'
$ws2.Cells.Item(18, 4).Value = 'Potential Biological Removal levels set as limit reference points, e.g. catch levels'
$ws2.Cells.Item(19, 1).Value = 22
$ws2.Cells.Item(19, 2).Value = 16
$ws2.Cells.Item(19, 3).Value = 'Catch Limit'
$ws2.Cells.Item(19, 4).Value = 'Two harvest control rules were evaluated in three uncertainty scenarios each.'
$ws2.Cells.Item(20, 1).Value = 23
$ws2.Cells.Item(20, 2).Value = 17
$ws2.Cells.Item(20, 3).Value = 'Catch Limit'
$ws2.Cells.Item(20, 4).Value = 'Harvest Control Rule with TAC adjusted based on stock assessment output, specified by zone or region (zone area > region area)'
$ws2.Cells.Item(21, 1).Value = 24
$ws2.Cells.Item(21, 2).Value = 18
$ws2.Cells.Item(21, 3).Value = 'Closure'
$ws2.Cells.Item(21, 4).Value = 'Spatial Closures, of three types.  Marine protected areas, species risk based closures, and triggered closures based on effort'
$ws2.Cells.Item(22, 1).Value = 25
$ws2.Cells.Item(22, 2).Value = 19
$ws2.Cells.Item(22, 3).Value = 'Catch Limit'
$ws2.Cells.Item(22, 4).Value = 'four Landing regulations:  landing obligation, 5% discard limit, year-to-year quota transfer, both'
$ws2.Cells.Item(23, 1).Value = 26
$ws2.Cells.Item(23, 2).Value = 20
$ws2.Cells.Item(23, 3).Value = 'Effort Limit'
$ws2.Cells.Item(23, 4).Value = 'None, the current method is used in all simulations'
$ws2.Cells.Item(24, 1).Value = 27
$ws2.Cells.Item(24, 2).Value = 21
$ws2.Cells.Item(24, 3).Value = 'Catch Limit'
$ws2.Cells.Item(24, 4).Value = 'Status quo, 2x, Profit max, Broken stick control rule, Spatial broken stick control rule'
$ws2.Cells.Item(25, 1).Value = 28
$ws2.Cells.Item(25, 2).Value = 21
$ws2.Cells.Item(25, 3).Value = 'Closure'
$ws2.Cells.Item(25, 4).Value = 'Spatial closures.  Closure trigger by zone or species with 20% or 30% triggers'
$ws2.Cells.Item(26, 1).Value = 29
$ws2.Cells.Item(26, 2).Value = 22
$ws2.Cells.Item(26, 3).Value = 'Catch Limit'
$ws2.Cells.Item(26, 4).Value = '8 approaches to setting P* buffer relative to F(lim).'
$ws2.Cells.Item(27, 1).Value = 30
$ws2.Cells.Item(27, 2).Value = 23
$ws2.Cells.Item(27, 3).Value = 'Catch Limit'
$ws2.Cells.Item(27, 4).Value = '5 harvest quotas options'
$ws2.Cells.Item(28, 1).Value = 31
$ws2.Cells.Item(28, 2).Value = 23
$ws2.Cells.Item(28, 3).Value = 'Closure'
$ws2.Cells.Item(28, 4).Value = 'Seasonal and area closures exist, but were constant across alternatives'
$ws2.Cells.Item(29, 1).Value = 32
$ws2.Cells.Item(29, 2).Value = 24
$ws2.Cells.Item(29, 3).Value = 'Catch Limit'
$ws2.Cells.Item(29, 4).Value = 'Catch limits set by HCR derived from monitoring methodology alternatives'
$ws2.Cells.Item(30, 1).Value = 33
$ws2.Cells.Item(30, 2).Value = 25
$ws2.Cells.Item(30, 3).Value = 'Catch Limit'
$ws2.Cells.Item(30, 4).Value = '10% reduction annually until recovery achieved or no change'
$ws2.Cells.Item(31, 1).Value = 34
$ws2.Cells.Item(31, 2).Value = 25
$ws2.Cells.Item(31, 3).Value = 'Effort Limit'
$ws2.Cells.Item(31, 4).Value = '10% reduction annually until recovery achieved or no change'
$ws2.Cells.Item(32, 1).Value = 35
$ws2.Cells.Item(32, 2).Value = 25
$ws2.Cells.Item(32, 3).Value = 'Closure'
$ws2.Cells.Item(32, 4).Value = 'Area closures, open or closed evaulated'
$ws2.Cells.Item(33, 1).Value = 36
$ws2.Cells.Item(33, 2).Value = 26
$ws2.Cells.Item(33, 3).Value = 'Catch Limit'
$ws2.Cells.Item(33, 4).Value = 'Limits set by a set of 4 Stock Assessment methodologies.'
$ws2.Cells.Item(34, 1).Value = 37
$ws2.Cells.Item(34, 2).Value = 27
$ws2.Cells.Item(34, 3).Value = 'Other'
$ws2.Cells.Item(34, 4).Value = 'Lamprey trapping options: Traps located in 14 or 16 streams, baited or unbaited.'
$ws2.Cells.Item(35, 1).Value = 39
$ws2.Cells.Item(35, 2).Value = 27
$ws2.Cells.Item(35, 3).Value = 'Other'
$ws2.Cells.Item(35, 4).Value = 'Lamprey baiting options: Used fixed bait amount or achieve a target in-water concentration.'
$ws2.Cells.Item(36, 1).Value = 40
$ws2.Cells.Item(36, 2).Value = 28
$ws2.Cells.Item(36, 3).Value = 'Closure;Size Limit'
$ws2.Cells.Item(36, 4).Value = 'Not evaluated by analysis, single alternative is a sex, size, and season regulation system'
$ws2.Cells.Item(37, 1).Value = 41
$ws2.Cells.Item(37, 2).Value = 29
$ws2.Cells.Item(37, 3).Value = 'Catch Limit'
$ws2.Cells.Item(37, 4).Value = 'Two TACC decision rule alternatives.'
$ws2.Cells.Item(38, 1).Value = 42
$ws2.Cells.Item(38, 2).Value = 30
$ws2.Cells.Item(38, 3).Value = 'Effort Limit'
$ws2.Cells.Item(38, 4).Value = '3 F based control rules.  Constant F, Reduced F when SSB<0.4B(0), Reduced F when SSB<0.7B(0),'
$ws2.Cells.Item(39, 1).Value = 43
$ws2.Cells.Item(39, 2).Value = 30
$ws2.Cells.Item(39, 3).Value = 'Effort Limit'
$ws2.Cells.Item(39, 4).Value = '7 baseline F levels.  0.1, 0.3, 0.5, 0.7, 1.0, 1.5, 2.0'
$ws2.Cells.Item(40, 1).Value = 44
$ws2.Cells.Item(40, 2).Value = 41
$ws2.Cells.Item(40, 3).Value = 'Access Control'
$ws2.Cells.Item(40, 4).Value = 'Constant across alternatives'
$ws2.Cells.Item(41, 1).Value = 45
$ws2.Cells.Item(41, 2).Value = 41
$ws2.Cells.Item(41, 3).Value = 'Closure'
$ws2.Cells.Item(41, 4).Value = '4 temporal closure alternatives.  Base case, add May closure, add October closure, closure rule - if sample below threshold catch close'
$ws2.Cells.Item(42, 1).Value = 46
$ws2.Cells.Item(42, 2).Value = 42
$ws2.Cells.Item(42, 3).Value = 'Catch Limit'
$ws2.Cells.Item(42, 4).Value = 'Two HCRs, based on a static or a dynamic B(0) used as a BRP.'
$ws2.Cells.Item(43, 1).Value = 47
$ws2.Cells.Item(43, 2).Value = 43
$ws2.Cells.Item(43, 3).Value = 'Catch Limit'
$ws2.Cells.Item(43, 4).Value = 'Two HCRs, one with a constant F, and one with a variable F based on BRPs'
$ws2.Cells.Item(44, 1).Value = 48
$ws2.Cells.Item(44, 2).Value = 43
$ws2.Cells.Item(44, 3).Value = 'Catch Limit'
$ws2.Cells.Item(44, 4).Value = 'Two additional HCRs, which correspond to the first two, but are adjusted based on environmental conditions.'
$ws2.Cells.Item(45, 1).Value = 49
$ws2.Cells.Item(45, 2).Value = 44
$ws2.Cells.Item(45, 3).Value = 'Catch Limit'
$ws2.Cells.Item(45, 4).Value = '4 management procedures covering different approaches to developing a rebuilding plan based on the acceptable recovery probability or inclusion of climate effects in the plan.'
$ws2.Cells.Item(46, 1).Value = 50
$ws2.Cells.Item(46, 2).Value = 45
$ws2.Cells.Item(46, 3).Value = 'Catch Limit'
$ws2.Cells.Item(46, 4).Value = 'This is a TAC based fishery.  The TACs weren''t the focus, but would be altered by reference point changes.'
$ws2.Cells.Item(47, 1).Value = 51
$ws2.Cells.Item(47, 2).Value = 45
$ws2.Cells.Item(47, 3).Value = 'Other'
$ws2.Cells.Item(47, 4).Value = 'Two reference point options, one that shifts to account for climate change, and one that does not.'
$ws2.Cells.Item(48, 1).Value = 52
$ws2.Cells.Item(48, 2).Value = 46
$ws2.Cells.Item(48, 3).Value = 'Catch Limit'
$ws2.Cells.Item(48, 4).Value = 'TAC implemented or not'
$ws2.Cells.Item(49, 1).Value = 53
$ws2.Cells.Item(49, 2).Value = 46
$ws2.Cells.Item(49, 3).Value = 'Size Limit'
$ws2.Cells.Item(49, 4).Value = 'Size limit implemented or not'
$ws2.Cells.Item(50, 1).Value = 54
$ws2.Cells.Item(50, 2).Value = 46
$ws2.Cells.Item(50, 3).Value = 'Closure'
$ws2.Cells.Item(50, 4).Value = 'marine protected areas implemented or not'
$ws2.Cells.Item(51, 1).Value = 55
$ws2.Cells.Item(51, 2).Value = 46
$ws2.Cells.Item(51, 3).Value = 'Other'
$ws2.Cells.Item(51, 4).Value = 'Pollution reduction implemented or not'
$ws2.Cells.Item(52, 1).Value = 56
$ws2.Cells.Item(52, 2).Value = 47
$ws2.Cells.Item(52, 3).Value = 'Other'
$ws2.Cells.Item(52, 4).Value = 'Note, not applicable to fishery management.  There are three approaches to beach replenishment, fixed amount and interval, fixed amount, and fixed interval, as well as no action.'
$ws2.Cells.Item(53, 1).Value = 57
$ws2.Cells.Item(53, 2).Value = 48
$ws2.Cells.Item(53, 3).Value = 'Closure'
$ws2.Cells.Item(53, 4).Value = 'Not cleared provided.'
$ws2.Cells.Item(54, 1).Value = 58
$ws2.Cells.Item(54, 2).Value = 31
$ws2.Cells.Item(54, 3).Value = 'Closure'
$ws2.Cells.Item(54, 4).Value = 'Spatial closures: No closures, 3 closure durations, 2 closure location rules'
$ws2.Cells.Item(55, 1).Value = 59
$ws2.Cells.Item(55, 2).Value = 31
$ws2.Cells.Item(55, 3).Value = 'Other'
$ws2.Cells.Item(55, 4).Value = 'Size based closure rules: 4 options'
$ws2.Cells.Item(56, 1).Value = 60
$ws2.Cells.Item(56, 2).Value = 49
$ws2.Cells.Item(56, 3).Value = 'Closure'
$ws2.Cells.Item(56, 4).Value = 'Spatial closures: No closures, 3 closure durations, 2 closure location rules'
$ws2.Cells.Item(57, 1).Value = 61
$ws2.Cells.Item(57, 2).Value = 49
$ws2.Cells.Item(57, 3).Value = 'Other'
$ws2.Cells.Item(57, 4).Value = 'Size based closure rules: 4 options'
$ws2.Cells.Item(58, 1).Value = 62
$ws2.Cells.Item(58, 2).Value = 50
$ws2.Cells.Item(58, 3).Value = 'Catch Limit'
$ws2.Cells.Item(58, 4).Value = '6 alternative methods for adjusting the catch limit.'
$ws2.Cells.Item(59, 1).Value = 63
$ws2.Cells.Item(59, 2).Value = 32
$ws2.Cells.Item(59, 3).Value = 'Access Control;Catch Limit'
$ws2.Cells.Item(59, 4).Value = 'Quota system, accompanied by location and gear restrictions, combined into 4 strategy'
$ws2.Cells.Item(60, 1).Value = 64
$ws2.Cells.Item(60, 2).Value = 33
$ws2.Cells.Item(60, 3).Value = 'Size Limit'
$ws2.Cells.Item(60, 4).Value = '# of size limits utilized and areas to which they are applied'
$ws2.Cells.Item(61, 1).Value = 65
$ws2.Cells.Item(61, 2).Value = 34
$ws2.Cells.Item(61, 3).Value = 'Catch Limit;Effort Limit'
$ws2.Cells.Item(61, 4).Value = '22 MPs available in the DLMtoolkit package in R, 11 output control MPs and 11 input control MPs'
$ws2.Cells.Item(62, 1).Value = 66
$ws2.Cells.Item(62, 2).Value = 35
$ws2.Cells.Item(62, 3).Value = 'Catch Limit'
$ws2.Cells.Item(62, 4).Value = 'TAC=ABC set based on ICES F based proceedures (with and without uncertainty) and a constant F'
$ws2.Cells.Item(63, 1).Value = 67
$ws2.Cells.Item(63, 2).Value = 36
$ws2.Cells.Item(63, 3).Value = 'Catch Limit'
$ws2.Cells.Item(63, 4).Value = '4 HCRs.  2 based on BRPs and 2 proportional harvest rules, the difference in each category is whether there is an annual TAC change limit.'
$ws2.Cells.Item(64, 1).Value = 68
$ws2.Cells.Item(64, 2).Value = 36
$ws2.Cells.Item(64, 3).Value = 'Other'
$ws2.Cells.Item(64, 4).Value = '3 stock assessment methods.  XSA, Schaefer, Difference'
$ws2.Cells.Item(65, 1).Value = 69
$ws2.Cells.Item(65, 2).Value = 37
$ws2.Cells.Item(65, 3).Value = 'Effort Limit'
$ws2.Cells.Item(65, 4).Value = '2 levels, Status quo and reduction from 9000 to 5000 boat days'
$ws2.Cells.Item(66, 1).Value = 70
$ws2.Cells.Item(66, 2).Value = 37
$ws2.Cells.Item(66, 3).Value = 'Closure'
$ws2.Cells.Item(66, 4).Value = '4 options, status quo, reef buffer area closure, Masig area closure, moon cycle calendar closure'
$ws2.Cells.Item(67, 1).Value = 71
$ws2.Cells.Item(67, 2).Value = 38
$ws2.Cells.Item(67, 3).Value = 'Catch Limit'
$ws2.Cells.Item(67, 4).Value = 'No alternatives, a single Harvest quota accompanied by a season duration'
$ws2.Cells.Item(68, 1).Value = 72
$ws2.Cells.Item(68, 2).Value = 39
$ws2.Cells.Item(68, 3).Value = 'Catch Limit'
$ws2.Cells.Item(68, 4).Value = '4 TACs; 750 t, 1000 t, 1250 t, and 1500 t; competitave TAC vs ITQ'
$ws2.Cells.Item(69, 1).Value = 73
$ws2.Cells.Item(69, 2).Value = 39
$ws2.Cells.Item(69, 3).Value = 'Effort Limit'
$ws2.Cells.Item(69, 4).Value = 'constant effort limit'
$ws2.Cells.Item(70, 1).Value = 74
$ws2.Cells.Item(70, 2).Value = 39
$ws2.Cells.Item(70, 3).Value = 'Closure'
$ws2.Cells.Item(70, 4).Value = 'four spatially explicit no-take extents: an extent consistent with that from the mid-1980s to mid-2004 (approximately 16% of coral trout habitat in the park); an extent implemented during rezoning in 2004 (32%); a hypothetical extent of 50% (Little et al. 2009a); and a hypothetical extent of 0%.'
$ws2.Cells.Item(71, 1).Value = 75
$ws2.Cells.Item(71, 2).Value = 40
$ws2.Cells.Item(71, 3).Value = 'Catch Limit'
$ws2.Cells.Item(71, 4).Value = '6 HCRs combining timeline and precaution:  a reactive decision interval with no additional ACL reduction, and five HCRs consisting of a fixed decision interval with precautionary ACL reductions of 0 (i.e., no reduction), 10, 20, 30, and 40%.'
$ws2.Cells.Item(72, 1).Value = 77
$ws2.Cells.Item(72, 2).Value = 51
$ws2.Cells.Item(72, 3).Value = 'Other'
$ws2.Cells.Item(72, 4).Value = '5 alternative methods to calculate biological reference points'
$ws2.Cells.Item(73, 1).Value = 78
$ws2.Cells.Item(73, 2).Value = 52
$ws2.Cells.Item(73, 3).Value = 'Other'
$ws2.Cells.Item(73, 4).Value = 'The primary alternatives evaluated here are different stock assessment methods.'
$ws2.Cells.Item(74, 1).Value = 79
$ws2.Cells.Item(74, 2).Value = 52
$ws2.Cells.Item(74, 3).Value = 'Catch Limit'
$ws2.Cells.Item(74, 4).Value = '3 alternative target escapement levels, a base case, as well as 20%, and 40% increases'
$ws2.Cells.Item(75, 1).Value = 80
$ws2.Cells.Item(75, 2).Value = 53
$ws2.Cells.Item(75, 3).Value = 'Catch Limit'
$ws2.Cells.Item(75, 4).Value = '5 Alternative harvest control rules (HCR), varrying by explotation rate and biological reference point time period'
$ws2.Cells.Item(76, 1).Value = 81
$ws2.Cells.Item(76, 2).Value = 54
$ws2.Cells.Item(76, 3).Value = 'Catch Limit'
$ws2.Cells.Item(76, 4).Value = 'A single HCR is evaluated for robust to alternative climate change scenarios'
$ws2.Cells.Item(77, 1).Value = 83
$ws2.Cells.Item(77, 2).Value = 55
$ws2.Cells.Item(77, 3).Value = 'Catch Limit'
$ws2.Cells.Item(77, 4).Value = 'A single harvest control rule is evaluated using alternative stock assessment models'

# Keep the active selection on the new sheet consistent with the source file.
$ws2.Range("C5").Select() | Out-Null

# Register the defined name that now points at the new sheet table range.
$wb.Names.Add("tblStudyManagementTools", "='tblStudyManagementTools1'!`$A`$1:`$D`$77")

Write-Output "Added sheet tblStudyManagementTools1 with $($ws2.UsedRange.Rows.Count) rows and defined name tblStudyManagementTools"
